# The commit swaps the content of ppt/theme/theme1.xml (the presentation's
# main theme, used by the slide master -> "Integral" / Red Violet) with
# ppt/theme/theme2.xml (only used by the notes master -> default "Office
# Theme"). Net visible effect: the slide master (and therefore every slide)
# switches from the "Integral" colour scheme to the standard Office colour
# scheme.
#
# The PowerPoint object model only exposes the deck's *active* theme
# through SlideMaster.Theme (ThemeColorScheme / ThemeFontScheme), so we
# reproduce the swap by writing the Office Theme's 12 theme colours onto
# that theme in place. The two themes already share an identical font
# scheme ("Office" / Arial) and format scheme, so only the colour scheme
# actually changes.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as PowerPoint RGB() integers (R + G*256 + B*65536).
$colorScheme.Colors(1).RGB  = 0         # dk1      #000000
$colorScheme.Colors(2).RGB  = 16777215  # lt1      #FFFFFF
$colorScheme.Colors(3).RGB  = 6968388   # dk2      #44546A
$colorScheme.Colors(4).RGB  = 15132391  # lt2      #E7E6E6
$colorScheme.Colors(5).RGB  = 13998939  # accent1  #5B9BD5
$colorScheme.Colors(6).RGB  = 3243501   # accent2  #ED7D31
$colorScheme.Colors(7).RGB  = 10855845  # accent3  #A5A5A5
$colorScheme.Colors(8).RGB  = 49407     # accent4  #FFC000
$colorScheme.Colors(9).RGB  = 12874308  # accent5  #4472C4
$colorScheme.Colors(10).RGB = 4697456   # accent6  #70AD47
$colorScheme.Colors(11).RGB = 12673797  # hlink    #0563C1
$colorScheme.Colors(12).RGB = 7491477   # folHlink #954F72

# Font scheme is already identical between the two themes (Arial-based
# "Office" font scheme), so nothing further to change there.
$fontScheme = $master.Theme.ThemeFontScheme
$fontScheme.MajorFont.Latin = "Arial"
$fontScheme.MinorFont.Latin = "Arial"
